$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("studyAmendments")
for ($r=1; $r -le 5; $r++) {
  $line = ""
  for ($c=1; $c -le 7; $c++) {
    $cell = $ws.Cells.Item($r, $c)
    $v = $cell.Value()
    $line += "[" + $cell.Address() + "=" + $v + "]"
  }
  Write-Output $line
}
